# Update "想去人数" (number of people interested) counts to the latest
# scraped values across the relevant worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1987
$wsExhibit.Range("F6").Value = 404
$wsExhibit.Range("F8").Value = 493
$wsExhibit.Range("F14").Value = 851
$wsExhibit.Range("F15").Value = 3776
$wsExhibit.Range("F17").Value = 824
$wsExhibit.Range("F18").Value = 601
$wsExhibit.Range("F21").Value = 1241
$wsExhibit.Range("F24").Value = 274

# Sheet "演出" (performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 42

# Sheet "全部类型" (all types, aggregated view)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 1987
$wsAll.Range("F10").Value = 404
$wsAll.Range("F12").Value = 493
$wsAll.Range("F18").Value = 851
$wsAll.Range("F21").Value = 3776
$wsAll.Range("F23").Value = 824
$wsAll.Range("F24").Value = 601
$wsAll.Range("F27").Value = 1241
$wsAll.Range("F30").Value = 274
